# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Replace the "Estado de Cuenta" detail rows (16-29) so that each
# worker's periods are grouped together (most recent period first),
# matching the refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; Doc = "9297632"; Nombre = "SAMIR YESID MARTELO MONTES"; Periodo = "2111"; Valor = 27861 },
    @{ Row = 17; Doc = "9297632"; Nombre = "SAMIR YESID MARTELO MONTES"; Periodo = "2110"; Valor = 36341 },
    @{ Row = 18; Doc = "9297632"; Nombre = "SAMIR YESID MARTELO MONTES"; Periodo = "2109"; Valor = 36341 },
    @{ Row = 19; Doc = "9297632"; Nombre = "SAMIR YESID MARTELO MONTES"; Periodo = "2108"; Valor = 36341 },
    @{ Row = 20; Doc = "9297632"; Nombre = "SAMIR YESID MARTELO MONTES"; Periodo = "2107"; Valor = 36341 },
    @{ Row = 21; Doc = "9297632"; Nombre = "SAMIR YESID MARTELO MONTES"; Periodo = "2106"; Valor = 36341 },
    @{ Row = 22; Doc = "9297632"; Nombre = "SAMIR YESID MARTELO MONTES"; Periodo = "2105"; Valor = 36341 },
    @{ Row = 23; Doc = "3805295"; Nombre = "HENRY DE JESUS GOMEZ PATERNINA"; Periodo = "2111"; Valor = 27861 },
    @{ Row = 24; Doc = "3805295"; Nombre = "HENRY DE JESUS GOMEZ PATERNINA"; Periodo = "2110"; Valor = 36341 },
    @{ Row = 25; Doc = "3805295"; Nombre = "HENRY DE JESUS GOMEZ PATERNINA"; Periodo = "2109"; Valor = 36341 },
    @{ Row = 26; Doc = "3805295"; Nombre = "HENRY DE JESUS GOMEZ PATERNINA"; Periodo = "2108"; Valor = 36341 },
    @{ Row = 27; Doc = "3805295"; Nombre = "HENRY DE JESUS GOMEZ PATERNINA"; Periodo = "2107"; Valor = 36341 },
    @{ Row = 28; Doc = "3805295"; Nombre = "HENRY DE JESUS GOMEZ PATERNINA"; Periodo = "2106"; Valor = 36341 },
    @{ Row = 29; Doc = "3805295"; Nombre = "HENRY DE JESUS GOMEZ PATERNINA"; Periodo = "2105"; Valor = 36341 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("C$n").Value = $r.Doc
    $ws.Range("D$n").Value = $r.Nombre
    $ws.Range("E$n").Value = $r.Periodo
    $ws.Range("F$n").Value = $r.Valor
}
